$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Insert three new columns (F, G, H) before the old "Storage Location" column,
# shifting the old F ("Storage Location") and G ("Cassetto") columns to I and J.
$ws.Range("F1:H1").EntireColumn.Insert()

# Headers for the three new columns
$ws.Range("F1").Value = "Total Reactions"
$ws.Range("G1").Value = "Reactions Used"
$ws.Range("H1").Value = "Reactions Available"

# Per-row Total Reactions / Reactions Used values
$totals = @{
    2 = 24; 3 = 24; 4 = 24;
    5 = 96; 6 = 96; 7 = 96; 8 = 96;
    9 = 24; 10 = 96;
    11 = 96; 12 = 96; 13 = 96; 14 = 96; 15 = 96;
    17 = 12; 18 = 12;
    24 = 300
}

for ($r = 2; $r -le 26; $r++) {
    if ($totals.ContainsKey($r)) {
        $ws.Cells.Item($r, 6).Value = $totals[$r]
    }
    if (($r -ge 2 -and $r -le 18) -or $r -eq 24) {
        $ws.Cells.Item($r, 7).Value = 0
    }
}

# Reactions Available = Total Reactions - Reactions Used, filled down as a shared formula
$ws.Range("H2").Formula = "=F2-G2"
$ws.Range("H3:H26").Formula = "=F3-G3"

Write-Host "done"
